# Weekly price-log update: insert the newest week's observation as a new
# row at the top of this variety's block (row 417), pushing the existing
# rows down by one (417->418, ..., 449->450).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(417).Insert()

$ws.Range("A417").Value = 4
$ws.Range("B417").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C417").Value = "Los Lagos"
$ws.Range("D417").Value = 45013
$ws.Range("E417").Value = 10
$ws.Range("F417").Value = 100112045
$ws.Range("G417").Value = "Zapallo"
$ws.Range("H417").Value = "Paine"
$ws.Range("I417").Value = "1a (cosecha)"
$ws.Range("J417").Value = 1200
$ws.Range("K417").Value = 500
$ws.Range("L417").Value = 500
$ws.Range("M417").Value = 500
$ws.Range("N417").Value = "$/kilo (volumen en unidades)"
$ws.Range("O417").Value = "Región de O'Higgins"
$ws.Range("P417").Value = 500
$ws.Range("Q417").Value = 1
$ws.Range("R417").Value = "Hortaliza"
